# Re-order the "Recorded By" (column G) value lists so that the exact
# token "System" is moved to the front of the comma-separated list
# (placed immediately after any leading lowercase "system" token(s), if present).
#
# Example:
#   "dnasr281@gmail.com, System"              -> "System, dnasr281@gmail.com"
#   "backup@backdoor.com, System"              -> "System, backup@backdoor.com"
#   "system, backup@backdoor.com, System"      -> "system, System, backup@backdoor.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Find the last used row in column G ("Recorded By" column).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # -4162 = xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if (-not ($val -is [string])) { continue }
    if ($val -notmatch "System") { continue }

    # Split on ", " preserving each trimmed token.
    $parts = $val -split ","
    for ($i = 0; $i -lt $parts.Count; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    $idx = -1
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($parts[$i].Equals("System")) {
            $idx = $i
            break
        }
    }

    if ($idx -lt 0) { continue }

    # Remove the "System" token from its current position.
    $rest = @()
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($i -ne $idx) {
            $rest += $parts[$i]
        }
    }

    # Determine insertion position: after any leading lowercase "system" tokens.
    $insertPos = 0
    while (($insertPos -lt $rest.Count) -and ($rest[$insertPos].Equals("system"))) {
        $insertPos++
    }

    # Rebuild the list with "System" inserted at $insertPos.
    $newParts = @()
    for ($i = 0; $i -lt $insertPos; $i++) {
        $newParts += $rest[$i]
    }
    $newParts += "System"
    for ($i = $insertPos; $i -lt $rest.Count; $i++) {
        $newParts += $rest[$i]
    }

    $newVal = $newParts -join ", "

    if ($newVal -ne $val) {
        $cell.Value2 = $newVal
    }
}
